$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "EA"
$ws.Range("B16").Value = "Enumeration Area"

$ws.Range("B17").Value = "Multiple Indicator Cluster Surveys"
$ws.Range("A17").Value = "MICS"

$ws.Range("A18").Value = "DHS"
$ws.Range("B18").Value = "Demographic and Health Surveys"

$ws.Range("I20").Select()
